$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: add a new "Carnegie Classification" data element under the
#     "Research Output" criteria (S/N 4), including the note describing the
#     classification buckets. Fill this row out BEFORE row 7 so the shared
#     strings "Carnegie Classification" / "carnegie_basic" are registered in
#     that order. ---
$ws.Range("C9").Value = "Carnegie Classification"
$ws.Range("D9").Value = 120
$ws.Range("E9").Value = "school"
$ws.Range("F9").Value = "carnegie_basic"

$ws.Range("G9").WrapText = $true
$ws.Range("G9").Value = "1- Doctoral Universities: Very High Research Activity`n2- Doctoral Universities: High Research Activity`n3- Special Focus Four-Year: Research Institution`n4- Doctoral/Professional Universities`n5- Master" + [char]0x2019 + "s Colleges & Universities (Larger, Medium, Small Programs)"

$ws.Rows(9).RowHeight = 75

# --- Row 7: continuation row for "Research Output" - same field/dev_category
#     as the new Carnegie Classification element above (reuses the shared
#     strings just created). ---
$ws.Range("D7").Value = 120
$ws.Range("E7").Value = "school"
$ws.Range("F7").Value = "carnegie_basic"

# --- Leave the filter-table selection on the filled block at the bottom of
#     the sheet. ---
[void]$ws.Range("B20:G26").Select()
